$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "319.18", "0.06860") are stored as text, matching the source data
# which uses inline strings throughout, and preserves formatting such as
# trailing zeros and the double-dot thousands/decimal style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.264.56'
$ws.Range("E2").Value = '  -2.27%  '
$ws.Range("D3").Value = '1.867.79'
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '319.18'
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.4377'
$ws.Range("E7").Value = '  -4.56%  '
$ws.Range("D8").Value = '0.3704'
$ws.Range("E8").Value = '  -2.96%  '
$ws.Range("D9").Value = '0.07512'
$ws.Range("E9").Value = '  -2.60%  '
$ws.Range("D10").Value = '0.9397'
$ws.Range("E10").Value = '  -3.72%  '
$ws.Range("D11").Value = '21.44'
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("D12").Value = '1.950.16'
$ws.Range("E12").Value = '  +2.62%  '
$ws.Range("D13").Value = '6.723'
$ws.Range("E13").Value = '  -3.07%  '
$ws.Range("D14").Value = '5.447'
$ws.Range("E14").Value = '  -3.47%  '
$ws.Range("D15").Value = '0.06860'
$ws.Range("E15").Value = '  -2.31%  '
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '82.32'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").Value = '0.000009082'
$ws.Range("E18").Value = '  -3.97%  '
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '15.96'
$ws.Range("E20").Value = '  -3.82%  '
$ws.Range("D21").Value = '28.266.42'
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("D22").Value = '5.129'
$ws.Range("E22").Value = '  -3.08%  '
$ws.Range("D23").Value = '10.79'
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("D24").Value = '2.124.93'
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("D25").Value = '2.026'
$ws.Range("E25").Value = '  -3.22%  '
$ws.Range("D26").Value = '154.95'
$ws.Range("E26").Value = '  -1.92%  '
$ws.Range("D27").Value = '18.43'
$ws.Range("E27").Value = '  -2.81%  '
$ws.Range("D28").Value = '5.319'
$ws.Range("E28").Value = '  -5.57%  '
$ws.Range("D29").Value = '113.86'
$ws.Range("E29").Value = '  -3.06%  '
$ws.Range("D30").Value = '1.728'
$ws.Range("E30").Value = '  -5.76%  '
$ws.Range("D31").Value = '0.09037'
$ws.Range("E31").Value = '  -2.28%  '
$ws.Range("D32").Value = '0.7985'
$ws.Range("E32").Value = '  -7.42%  '
$ws.Range("D33").Value = '4.834'
$ws.Range("E33").Value = '  -4.86%  '
$ws.Range("D34").Value = '1.171'
$ws.Range("E34").Value = '  -5.58%  '
$ws.Range("D35").Value = '2.957'
$ws.Range("E35").Value = '  -1.30%  '
$ws.Range("D36").Value = '1.002'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").Value = '1.119'
$ws.Range("E37").Value = '  -2.43%  '
$ws.Range("D38").Value = '0.05431'
$ws.Range("E38").Value = '  -4.79%  '
$ws.Range("D39").Value = '0.01955'
$ws.Range("E39").Value = '  -3.59%  '
$ws.Range("D40").Value = '2.962'
$ws.Range("E40").Value = '  +7.30%  '
$ws.Range("D41").Value = '7.129'
$ws.Range("E41").Value = '  -3.45%  '
$ws.Range("D42").Value = '0.5253'
$ws.Range("E42").Value = '  -4.15%  '
$ws.Range("D43").Value = '0.1672'
$ws.Range("E43").Value = '  -4.45%  '
$ws.Range("D44").Value = '8.712'
$ws.Range("E44").Value = '  -5.93%  '
$ws.Range("D45").Value = '0.06765'
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("D46").Value = '2.040'
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("D47").Value = '0.4864'
$ws.Range("E47").Value = '  -5.63%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '10.58'
$ws.Range("E48").Value = '  -6.17%  '
$ws.Range("D49").Value = '107.91'
$ws.Range("E49").Value = '  -1.99%  '
$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = '0.000002502'
$ws.Range("E50").Value = '  -4.77%  '
$ws.Range("D51").Value = '1.678'
$ws.Range("E51").Value = '  -5.20%  '
